$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 232, shifting rows 232:318 down to 233:319
$ws.Rows("232:232").Insert()

# Fill in the new row 232 with data (values copied from the row above for
# the constant columns, new values for the variable columns)
$ws.Range("A232").Value = 4
$ws.Range("B232").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C232").Value = "Los Lagos"
$ws.Range("D232").Value = 44900
$ws.Range("E232").Value = 10
$ws.Range("F232").Value = 100112044
$ws.Range("G232").Value = "Perejil"
$ws.Range("H232").Value = "Sin especificar"
$ws.Range("I232").Value = "Primera"
$ws.Range("J232").Value = 70
$ws.Range("K232").Value = 6000
$ws.Range("L232").Value = 6000
$ws.Range("M232").Value = 6000
$ws.Range("N232").Value = "`$/docena de atados (2 kilos)"
$ws.Range("O232").Value = "Región de La Araucanía"
$ws.Range("P232").Value = 3000
$ws.Range("Q232").Value = 2
$ws.Range("R232").Value = "Hortaliza"

# Match the date-cell style used throughout column D
$ws.Range("D232").NumberFormat = $ws.Range("D233").NumberFormat
